# [Feat 2269] Add support of steps worksheet metadata.
# Rename the "STEPS TODO" sheet to "STEPS" and populate its header row
# with the step-related column names (replacing the former string
# literals used during export with the matching enum names), then make
# it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STEPS TODO")
$ws.Name = "STEPS"

# Populate the header row in an order that mirrors the shared-string
# insertion order of the original authored workbook.
$ws.Range("A1").Value = "ACTION"
$ws.Range("D1").Value = "TC_STEP_ID"
$ws.Range("E1").Value = "TC_STEP_NUM"
$ws.Range("G1").Value = "TC_STEP_ACTION"
$ws.Range("I1").Value = "TC_STEP_#_REQ"
$ws.Range("J1").Value = "TC_STEP_#_ATTACHMENT"
$ws.Range("K1").Value = "TC_STEP_CUF_<CODE>"
$ws.Range("B1").Value = "TC_OWNER_PATH"
$ws.Range("C1").Value = "TC_OWNER_ID"
$ws.Range("F1").Value = "TC_STEP_IS_CALL_STEP"
$ws.Range("H1").Value = "TC_STEP_EXPECTED_RESULT"

# Best-fit the header columns to their content (as Excel would after
# typing the headers and auto-sizing the columns).
$bestFitWidths = @(7.0, 15.666666666666666, 12.666666666666666, 10.333333333333334, 12.833333333333334, 20.166666666666668, 15.5, 24.5, 14.0, 22.833333333333332, 19.833333333333332)
for ($i = 0; $i -lt $bestFitWidths.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).ColumnWidth = $bestFitWidths[$i]
}

# Make STEPS the active sheet/tab (moves tabSelected from TEST_CASES to
# STEPS and records activeTab on the workbook view).
$ws.Activate()
$ws.Range("K18").Select() | Out-Null
